# Edit script: updates the Victoria key outbreaks cluster table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("A1").Value = "Cluster"
$ws.Range("B1").Value = "Active cases"

# Data: row -> (Cluster name, Active cases)
$data = @(
    @("3153 Sacred Heart Community St Kilda Tier 1A", 21),
    @("3600 Belvedere Aged Care Noble Park", 39),
    @("3612 BlueCross Glengowrie", 56),
    @("3684 Homestyle Aged Care Langford Grange Cranbourne East", 30),
    @("3824 Estia Health South Morang", 32),
    @("3980 Arcare Keysborough Aged Care Keysborough", 27),
    @("4075 Homestyle Aged Care Ferndale Gardens Aged Care Services Bayswater North", 16),
    @("4518 Regis Aged Care Fawkner", 21),
    @("ACFS Port Logistics Altona", 14),
    @("Bespoke Childcare Dingley Village", 13),
    @("Community Kids Pascoe Vale Early Education Centre Pascoe Vale", 23),
    @("Guardian Childcare Caulfield", 17),
    @("Hello Fresh Warehouse Ravenhall", 18),
    @("Honeyeater Hairdressers Bendigo", 10),
    @("Inghams Enterprises Somerville", 27),
    @("JBS Australia Brooklyn", 10),
    @("Metro Tunnel Shine Domain Site Albert Road Construction Site South Melbourne", 10),
    @("Northern Health The Northern Hospital Epping", 17),
    @("Shawlands Caravan Park Dandenong South", 17),
    @("St Vincents Hospital Emergency Department Melbourne", 32),
    @("The Robin Hood Inn Drouin West", 49),
    @("The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1A", 10),
    @("The Royal Children's Hospital Parkville", 11),
    @("Visy Recycling Springvale", 10),
    @("Werribee Mercy Hospital Emergency Department", 33),
    @("Western Health Sunshine Hospital Emergency Department", 22)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# Original table had 29 rows; new table only has 27. Delete the now-removed rows 28 and 29.
$lastOldRow = 29
$lastNewRow = $row - 1
if ($lastOldRow -gt $lastNewRow) {
    $deleteRange = $ws.Range($ws.Cells.Item($lastNewRow + 1, 1), $ws.Cells.Item($lastOldRow, 2))
    $deleteRange.Delete()
}
